# code fix for support
# Applies the Portfolio / Faculty / V2Project / Hsr data-entry changes plus
# the resulting sheet selections, matching the recorded Excel session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Portfolio sheet: fill in row 2 with a freshly-generated automation
# portfolio record (previously only AccountName was present).
# ---------------------------------------------------------------------
$wsPortfolio = $wb.Worksheets.Item("Portfolio")
$wsPortfolio.Activate()
$wsPortfolio.Range("A2").Value  = "Automation portfolio 4377747"
$wsPortfolio.Range("B2").Value  = "Government"
$wsPortfolio.Range("C2").Value  = "100 RESILIENT CITIES"
$wsPortfolio.Range("D2").Value  = "15"
$wsPortfolio.Range("E2").Value  = "10"
$wsPortfolio.Range("F2").Value  = "26556"
$wsPortfolio.Range("G2").Value  = "United States"
$wsPortfolio.Range("H2").Value  = "Ohio"
$wsPortfolio.Range("I2").Value  = "88741 Lucas Locks"
$wsPortfolio.Range("J2").Value  = "Durganberg"
$wsPortfolio.Range("K2").Value  = "59989-8155"
$wsPortfolio.Range("L2").Value  = "Joey"
$wsPortfolio.Range("M2").Value  = "jayne.keebler@hotmail.com"
$wsPortfolio.Range("N2").Value  = "7362255383"
$wsPortfolio.Range("C1").Select()

# ---------------------------------------------------------------------
# Faculty sheet: add the trailing "Org" column (I) with value "Testing".
# ---------------------------------------------------------------------
$wsFaculty = $wb.Worksheets.Item("Faculty")
$wsFaculty.Activate()
$wsFaculty.Range("I1").Value = "Org"
$wsFaculty.Range("I2").Value = "Testing"
$wsFaculty.Range("A4").Select()

# ---------------------------------------------------------------------
# Hsr sheet: add the trailing "Org" column (F) header only.
# ---------------------------------------------------------------------
$wsHsr = $wb.Worksheets.Item("Hsr")
$wsHsr.Activate()
$wsHsr.Range("F1").Value = "Org"
$wsHsr.Range("F1").Select()

# ---------------------------------------------------------------------
# V2Project sheet: regenerate ProjectId, append FeatureName/Subject
# columns (T/U) with the new ticket's subject line.
# ---------------------------------------------------------------------
$wsV2 = $wb.Worksheets.Item("V2Project")
$wsV2.Activate()
$wsV2.Range("A2").NumberFormat = "@"
$wsV2.Range("A2").Value = "2202266502"
$wsV2.Range("T1").Value = "FeatureName"
$wsV2.Range("U1").Value = "Subject"
$wsV2.Range("T2").Value = "A01 Air Quality"
$wsV2.Range("E8").Select()
